# Apply the "notes for reflection report.docx" edits:
#  1) merge "Difficulty getting <spellcheck>box_num</spellcheck> calculations..."
#     into a single run (drops the proofErr spell-check wrapper)
#  2) merge "Created <spellcheck>partially_complete</spellcheck> method..." into
#     a single run
#  3) merge the quoted "create_box_num_list" span into a single run, and split
#     the following sentence ("If loop for moving pointer...") into several
#     short runs, within the "Considered how to generate..." paragraph
#  4) merge "After creating <spellcheck>create_box_num_list</spellcheck> ...
#     <spellcheck>check_valid</spellcheck> method..." into a single run
#  5) append five new reflection paragraphs (plus spacer blank paragraphs)
#     describing the candidate-structure / find_match work

$d = $word.ActiveDocument

function Get-ParagraphStartingWith($prefix) {
    foreach ($para in $d.Paragraphs) {
        if ($para.Range.Text.StartsWith($prefix)) {
            return $para
        }
    }
    return $null
}

# --- Change 1: "Difficulty getting box_num calculations..." paragraph ---
# Searching for the full (post-edit) text and "replacing" it with itself makes
# Word re-flow the matched range as one run, which also drops the proofErr
# spell-check markers that used to bracket "box_num".
$old1 = "Difficulty getting box_num calculations to work, eventually realised needed to do calculation as int to truncate decimals"
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $old1, 2) | Out-Null

# --- Change 2: "Created partially_complete method..." paragraph ---
$old2 = "Created partially_complete method, managed to fill correct boxes with same int for all sizes of board. "
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $old2, 2) | Out-Null

# --- Changes 3 & 4: "Considered how to generate..." paragraph ---
# This paragraph needs (a) the quoted "create_box_num_list" span merged into
# one run (dropping its proofErr wrapper), and (b) the trailing sentence split
# into several new runs. Replace the whole paragraph body via InsertXML so the
# resulting run layout is exact.
$pConsidered = Get-ParagraphStartingWith("Considered how to generate")
$pConsidered.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Considered how to generate unique random numbers to fill each box, created doubly linked</w:t></w:r><w:r><w:t xml:space="preserve"> circular</w:t></w:r><w:r><w:t xml:space="preserve"> list </w:t></w:r><w:r><w:t>“create_box_num_list” method to store all numbers with idea to remove each number as selected</w:t></w:r><w:r><w:t xml:space="preserve">. Used circular so do not need to worry about </w:t></w:r><w:r><w:t xml:space="preserve">where the list pointer is pointing, and just move the pointer along </w:t></w:r><w:r><w:t xml:space="preserve">in list by </w:t></w:r><w:r><w:t>a random number to max number of elements in</w:t></w:r><w:r><w:t xml:space="preserve"> list -1 (</w:t></w:r><w:r><w:t>no point looping through to point back at the same element)</w:t></w:r><w:r><w:t>. If</w:t></w:r><w:r><w:t>-</w:t></w:r><w:r><w:t>loop for moving pointer</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> only moves pointer when there are greater than 1 element remaining so not to waste time by moving pointer to point at itself.</w:t></w:r></w:p>') | Out-Null

# --- Change 5: "After creating create_box_num_list function..." paragraph ---
$old5 = "After creating create_box_num_list function and node structure, realised could create a structure to hold all the candidates for rows, columns and boxes. This could then be used to complete / solve the board. It would also remove the requirement for a check_valid method, as everything inserted would have to be valid if its still a candidate (just look for matches between the row/box/column candidates). It would also reduce the time to compute as now not just trying any number, but only candidates."
$d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $old5, 2) | Out-Null

# --- Change 6: add the new reflection paragraphs after that paragraph ---
$pAfterCreating = Get-ParagraphStartingWith("After creating create_box_num_list")
$insertPos = $pAfterCreating.Range.End
$insertRange = $d.Range($insertPos, $insertPos)
$insertRange.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Function created to populate candidate structures for each box/row/column</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Rewrote partially_complete method to use the candidate structures and remove candidates as entering each into grid.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:lastRenderedPageBreak/><w:t>Was using circular doubly linked list so didn’t need to worry where pointer was in list during shuffling for population of non conflicting boxes on the board. However when find_match function, realised this meant it wasn’t possible to check if number in list was greater than number comparing to, as could already have started ahead of it. This would lead to more comparisons than necessary. So rewrote functions so candidates list was now a non circular doubly linked list. The candidates structure will now always be pointing at the smallest element in the list.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Wrote find_match function, originally was calling set_order_to_compare to find shortest, mid and longest lists from within function. However realised this will make it difficult to find ALL matches later on when checking number of solutions. Now set_order_to_compare will be called externally to find_match, and the lists can be passed to the find_match function. This will allow the pointer to the shortest list to be advanced before passing so as to be able to utilise the function to find ALL matches.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>') | Out-Null
